$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row stays the same (codice_1_livello, label _ITA _1 _livello, label_ENG_1_livello, definizione_ITA, definizione_ENG) ---
# (no change needed to row 1)

# --- Row 2: Mrs / Sig.ra ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Sig.ra"
$ws.Range("C2").Value = "Mrs"
$ws.Range("D2").Value = "Abbreviazione di 'Signora'."
$ws.Range("E2").Value = "English honorific used for women, usually for those who are married and who do not instead use another title."

# --- Row 3: Miss / Sig.na ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Sig.na"
$ws.Range("C3").Value = "Miss"
$ws.Range("D3").Value = "Abbreviazione di 'Signorina'."
$ws.Range("E3").Value = "English honorific used for not married women."

# --- Row 4: Mr / Sig ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Sig"
$ws.Range("C4").Value = "Mr"
$ws.Range("D4").Value = "Abbreviazione di 'Signore'."
$ws.Range("E4").Value = "Abbreviated form of 'Mister'."

# --- Row 5: Ms (new item, no Italian equivalent) ---
$ws.Range("A5").Value = 4
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value = "Ms"
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = "English honorific used with the last name or full name of a woman, intended as a default form of address for women regardless of marital status"

# --- Row 6: Dr / Dott ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Dott"
$ws.Range("C6").Value = "Dr"
$ws.Range("D6").Value = "Abbreviazione di 'Dottore'."
$ws.Range("E6").Value = "Abbreviated form of 'Doctor'."

# --- Row 7: Dott.ssa (Italian-only item) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Dott.ssa"
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = "Abbreviazione di 'Dottoressa'."
$ws.Range("E7").ClearContents()

# --- Row 8: Prof / Prof (new row, moved down) ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Prof"
$ws.Range("C8").Value = "Prof"
$ws.Range("D8").Value = "Abbreviazione di 'Professore'."
$ws.Range("E8").Value = "Abbreviated form of 'Professor'."

# --- Row 9: Prof.ssa (new item, Italian-only) ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Prof.ssa"
$ws.Range("D9").Value = "Abbreviazione di 'Professoressa'"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 18.333333333333332
$ws.Columns.Item(3).ColumnWidth = 18.333333333333332
$ws.Columns.Item(4).ColumnWidth = 28.666666666666668

# --- View: zoom + selection ---
$excel.ActiveWindow.Zoom = 194
$ws.Range("D10").Select()
